# Applies the scheduled-runner price/profit data refresh captured in the
# commit "chore: update Sheets via scheduled runner" across all 8 Disciple
# of the Hand profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each row's H:N columns (currentAveragePrice.. / LevePrice.. / LeveProfit..)
# are refreshed with newly-sampled market data; a few rows gain/lose their
# HQ profit (column N) cell depending on whether an HQ price was available.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
# row 96
$wsALC.Range("H96").Value = 1098.6666
$wsALC.Range("I96").Value = 1138.4
$wsALC.Range("K96").Value = 3415.2
$wsALC.Range("M96").Value = -2042.2
# row 100
$wsALC.Range("H100").Value = 205661.2
$wsALC.Range("I100").Value = 336100
$wsALC.Range("J100").Value = 10003
$wsALC.Range("K100").Value = 336100
$wsALC.Range("L100").Value = 10003
$wsALC.Range("M100").Value = -335559
$wsALC.Range("N100").Value = -11085
# row 111
$wsALC.Range("H111").Value = 3507.6667
$wsALC.Range("I111").Value = 3747.5
$wsALC.Range("J111").Value = 3028
$wsALC.Range("K111").Value = 11242.5
$wsALC.Range("L111").Value = 9084
$wsALC.Range("M111").Value = -8175.5
$wsALC.Range("N111").Value = -15218
# row 113
$wsALC.Range("H113").Value = 6920
$wsALC.Range("I113").Value = 6494.4287
$wsALC.Range("J113").Value = 7416.5
$wsALC.Range("K113").Value = 6494.4287
$wsALC.Range("L113").Value = 7416.5
$wsALC.Range("M113").Value = -3240.4287
$wsALC.Range("N113").Value = -13924.5
# row 116
$wsALC.Range("H116").Value = 6012.8667
$wsALC.Range("J116").Value = 6470.643
$wsALC.Range("L116").Value = 6470.643
$wsALC.Range("N116").Value = -13354.643
# row 132
$wsALC.Range("H132").Value = 1911.7906
$wsALC.Range("I132").Value = 1961
$wsALC.Range("J132").Value = 1537.8
$wsALC.Range("K132").Value = 5883
$wsALC.Range("L132").Value = 4613.4
$wsALC.Range("M132").Value = -3353
$wsALC.Range("N132").Value = -9673.4

# --- ARM ---
# row 32
$wsARM.Range("H32").Value = 3648.1594
$wsARM.Range("I32").Value = 3648.1594
$wsARM.Range("K32").Value = 3648.1594
$wsARM.Range("M32").Value = -3361.1594
# row 102
$wsARM.Range("H102").Value = 2281.5
$wsARM.Range("I102").Value = 2281.5
$wsARM.Range("K102").Value = 2281.5
$wsARM.Range("M102").Value = -659.5
# row 132
$wsARM.Range("H132").Value = 5138.08
$wsARM.Range("I132").Value = 3026.524
$wsARM.Range("K132").Value = 9079.572
$wsARM.Range("M132").Value = -6549.572

# --- BSM ---
# row 19
$wsBSM.Range("H19").Value = 3416.6667
# row 80
$wsBSM.Range("H80").Value = 264.2
$wsBSM.Range("I80").Value = 235.85715
$wsBSM.Range("J80").Value = 289
$wsBSM.Range("K80").Value = 235.85715
$wsBSM.Range("L80").Value = 289
$wsBSM.Range("M80").Value = 762.14285
$wsBSM.Range("N80").Value = -2285
# row 83
$wsBSM.Range("H83").Value = 264.2
$wsBSM.Range("I83").Value = 235.85715
$wsBSM.Range("J83").Value = 289
$wsBSM.Range("K83").Value = 1179.28575
$wsBSM.Range("L83").Value = 1445
$wsBSM.Range("M83").Value = 3812.71425
$wsBSM.Range("N83").Value = -11429
# row 99
$wsBSM.Range("H99").Value = 5568.0625
$wsBSM.Range("I99").Value = 5568.0625
$wsBSM.Range("K99").Value = 5568.0625
$wsBSM.Range("M99").Value = -4070.0625
# row 107
$wsBSM.Range("H107").Value = 4466.143
$wsBSM.Range("I107").Value = 2321.9092
$wsBSM.Range("J107").Value = 6824.8
$wsBSM.Range("K107").Value = 2321.9092
$wsBSM.Range("L107").Value = 6824.8
$wsBSM.Range("M107").Value = -401.9092000000001
$wsBSM.Range("N107").Value = -10664.8
# row 134
$wsBSM.Range("H134").Value = 8150.469
$wsBSM.Range("I134").Value = 3190.8235
$wsBSM.Range("K134").Value = 9572.470499999999
$wsBSM.Range("M134").Value = -7037.470499999999

# --- CRP ---
# row 58
$wsCRP.Range("H58").Value = 3156
$wsCRP.Range("J58").Value = 4442.625
$wsCRP.Range("L58").Value = 4442.625
$wsCRP.Range("N58").Value = -4848.625
# row 99
$wsCRP.Range("H99").Value = 8871.625
$wsCRP.Range("I99").Value = 4776.8887
$wsCRP.Range("K99").Value = 4776.8887
$wsCRP.Range("M99").Value = -3278.8887
# row 105
$wsCRP.Range("H105").Value = 1756.5555
$wsCRP.Range("I105").Value = 1837
$wsCRP.Range("K105").Value = 1837
$wsCRP.Range("M105").Value = -90
# row 126
$wsCRP.Range("H126").Value = 8871.625
$wsCRP.Range("I126").Value = 4776.8887
$wsCRP.Range("K126").Value = 14330.6661
$wsCRP.Range("M126").Value = -11860.6661
# row 134
$wsCRP.Range("H134").Value = 3629.9773
$wsCRP.Range("I134").Value = 2385.1614
$wsCRP.Range("J134").Value = 6598.385
$wsCRP.Range("K134").Value = 7155.4842
$wsCRP.Range("L134").Value = 19795.155
$wsCRP.Range("M134").Value = -4620.4842
$wsCRP.Range("N134").Value = -24865.155
# row 136
$wsCRP.Range("H136").Value = 3156
$wsCRP.Range("J136").Value = 4442.625
$wsCRP.Range("L136").Value = 13327.875
$wsCRP.Range("N136").Value = -18427.875

# --- CUL ---
# row 122
$wsCUL.Range("H122").Value = 4167325
$wsCUL.Range("I122").Value = 624.3333
$wsCUL.Range("J122").Value = 8334026
$wsCUL.Range("K122").Value = 5618.9997
$wsCUL.Range("L122").Value = 75006234
$wsCUL.Range("M122").Value = -3168.9997
$wsCUL.Range("N122").Value = -75011134
# row 131
$wsCUL.Range("H131").Value = 3131.963
$wsCUL.Range("J131").Value = 4271.353
$wsCUL.Range("L131").Value = 12814.059
$wsCUL.Range("N131").Value = -22894.059

# --- GSM ---
# row 39
$wsGSM.Range("H39").Value = 53149.332
$wsGSM.Range("J39").Value = 53149.332
$wsGSM.Range("L39").Value = 53149.332
$wsGSM.Range("N39").Value = -54213.332
# row 69
$wsGSM.Range("H69").Value = 0
$wsGSM.Range("J69").Value = 0
$wsGSM.Range("L69").Value = 0
$wsGSM.Range("N69").ClearContents()
# row 72
$wsGSM.Range("H72").Value = 0
$wsGSM.Range("J72").Value = 0
$wsGSM.Range("L72").Value = 0
$wsGSM.Range("N72").ClearContents()
# row 80
$wsGSM.Range("H80").Value = 2599.8
$wsGSM.Range("I80").Value = 2599.6667
$wsGSM.Range("J80").Value = 2600
$wsGSM.Range("K80").Value = 2599.6667
$wsGSM.Range("L80").Value = 2600
$wsGSM.Range("M80").Value = -1601.6667
$wsGSM.Range("N80").Value = -4596
# row 83
$wsGSM.Range("H83").Value = 2599.8
$wsGSM.Range("I83").Value = 2599.6667
$wsGSM.Range("J83").Value = 2600
$wsGSM.Range("K83").Value = 12998.3335
$wsGSM.Range("L83").Value = 13000
$wsGSM.Range("M83").Value = -8006.333500000001
$wsGSM.Range("N83").Value = -22984
# row 97
$wsGSM.Range("H97").Value = 1029.75
$wsGSM.Range("I97").Value = 770
$wsGSM.Range("K97").Value = 770
$wsGSM.Range("M97").Value = -274
# row 102
$wsGSM.Range("H102").Value = 2977.5
$wsGSM.Range("I102").Value = 2977.5
$wsGSM.Range("K102").Value = 2977.5
$wsGSM.Range("M102").Value = -1355.5
# row 113
$wsGSM.Range("H113").Value = 102906.7
$wsGSM.Range("I113").Value = 136172.6
$wsGSM.Range("J113").Value = 3109
$wsGSM.Range("K113").Value = 136172.6
$wsGSM.Range("L113").Value = 3109
$wsGSM.Range("M113").Value = -134002.6
$wsGSM.Range("N113").Value = -7449
# row 120
$wsGSM.Range("H120").Value = 64999
$wsGSM.Range("J120").Value = 64999
$wsGSM.Range("L120").Value = 64999
$wsGSM.Range("N120").Value = -74675
# row 122
$wsGSM.Range("H122").Value = 3665.3333
$wsGSM.Range("I122").Value = 3748.5
$wsGSM.Range("K122").Value = 11245.5
$wsGSM.Range("M122").Value = -8795.5
# row 126
$wsGSM.Range("H126").Value = 5749.5
$wsGSM.Range("I126").Value = 1500
$wsGSM.Range("J126").Value = 9999
$wsGSM.Range("K126").Value = 4500
$wsGSM.Range("L126").Value = 29997
$wsGSM.Range("M126").Value = -2030
$wsGSM.Range("N126").Value = -34937

# --- LTW ---
# row 132
$wsLTW.Range("H132").Value = 7647.1875
$wsLTW.Range("I132").Value = 7311.7856
$wsLTW.Range("K132").Value = 21935.3568
$wsLTW.Range("M132").Value = -19405.3568
# row 136
$wsLTW.Range("H136").Value = 4175.8604
$wsLTW.Range("J136").Value = 6741.875
$wsLTW.Range("L136").Value = 20225.625
$wsLTW.Range("N136").Value = -25325.625

# --- WVR ---
# row 9
$wsWVR.Range("H9").Value = 12164.5
$wsWVR.Range("I9").Value = 13886.333
$wsWVR.Range("K9").Value = 13886.333
$wsWVR.Range("M9").Value = -13746.333
# row 107
$wsWVR.Range("H107").Value = 1401.826
$wsWVR.Range("I107").Value = 1060.3889
$wsWVR.Range("K107").Value = 3181.1667
$wsWVR.Range("M107").Value = -1261.1667
# row 126
$wsWVR.Range("H126").Value = 9472.3125
$wsWVR.Range("J126").Value = 17319.777
$wsWVR.Range("L126").Value = 51959.33099999999
$wsWVR.Range("N126").Value = -56899.33099999999
# row 131
$wsWVR.Range("H131").Value = 46571.668
$wsWVR.Range("J131").Value = 46571.668
$wsWVR.Range("L131").Value = 46571.668
$wsWVR.Range("N131").Value = -56651.668
# row 132
$wsWVR.Range("H132").Value = 124312.16
$wsWVR.Range("I132").Value = 196738.27
$wsWVR.Range("J132").Value = 20566.648
$wsWVR.Range("K132").Value = 590214.8099999999
$wsWVR.Range("L132").Value = 61699.944
$wsWVR.Range("M132").Value = -587684.8099999999
$wsWVR.Range("N132").Value = -66759.944
# row 136
$wsWVR.Range("H136").Value = 7410030
$wsWVR.Range("I136").Value = 11113584
$wsWVR.Range("J136").Value = 2921.3333
$wsWVR.Range("K136").Value = 33340752
$wsWVR.Range("L136").Value = 8763.999899999999
$wsWVR.Range("M136").Value = -33338202
$wsWVR.Range("N136").Value = -13863.9999
